$wb = $excel.ActiveWorkbook

# ===================== Sheet 1: Overview =====================
$ws1 = $wb.Worksheets.Item("Overview")

# Clear existing hyperlinks on the sheet (Range.Hyperlinks.Delete removes all on sheet)
$ws1.Range("A1").Hyperlinks.Delete()

# Update non-hyperlink status cells (B/C columns)
$ws1.Range("B2").Value = "Handed back: in sync with en-US"
$ws1.Range("C2").Value = "Handed back: in sync with en-US"
$ws1.Range("B3").Value = "Handed back: in sync with en-US"
$ws1.Range("C3").Value = "Handed back: in sync with en-US"
$ws1.Range("B4").Value = "Ready for handoff"
$ws1.Range("C4").Value = "Ready for handoff"
$ws1.Range("B5").Value = "Not to be localized"
$ws1.Range("C5").Value = "Not to be localized"

# Re-create hyperlinks with updated display text, same underlying targets
$ws1.Hyperlinks.Add($ws1.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/11886f77cf3d268b498142352f706872679342af/e2e/c3f6fd5a-737b-46c9-98c1-5f206e633e5e.md", "", "", "ffff4f896d8c-d0bd-4be3-aacb-602751cc319c.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/11886f77cf3d268b498142352f706872679342af/e2e/ffff4f896d8c-d0bd-4be3-aacb-602751cc319c.md", "", "", "ffffffc9b3ff4e-2f02-46c4-a776-5467936972a9.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/11886f77cf3d268b498142352f706872679342af/e2e/ffffffc9b3ff4e-2f02-46c4-a776-5467936972a9.md", "", "", "c3f6fd5a-737b-46c9-98c1-5f206e633e5e.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/11886f77cf3d268b498142352f706872679342af/.localization-config", "", "", ".localization-config") | Out-Null

# ===================== Sheet 2: zh-cn =====================
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Range("A1").Hyperlinks.Delete()

$ws2.Range("B2").Value = "Handed back: in sync with en-US"
$ws2.Range("D2").Value = "2016-01-26 12:17:13"
$ws2.Range("G2").Value = "2016-01-26 12:18:00"
$ws2.Range("H2").Value = "Include"
$ws2.Range("B3").Value = "Handed back: in sync with en-US"
$ws2.Range("D3").Value = "2016-01-26 12:17:13"
$ws2.Range("G3").Value = "2016-01-26 12:18:00"
$ws2.Range("H3").Value = "Include"
$ws2.Range("B4").Value = "Ready for handoff"
$ws2.Range("D4").Value = "2016-01-26 12:23:23"
$ws2.Range("G4").Value = "2016-01-26 12:22:08"
$ws2.Range("H4").Value = "Include"
$ws2.Range("B5").Value = "Not to be localized"
$ws2.Range("D5").Value = "0001-01-01 00:00:00"
$ws2.Range("G5").Value = "0001-01-01 00:00:00"
$ws2.Range("H5").Value = "Ignored"

$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/11886f77cf3d268b498142352f706872679342af/e2e/c3f6fd5a-737b-46c9-98c1-5f206e633e5e.md", "", "", "ffff4f896d8c-d0bd-4be3-aacb-602751cc319c.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8b187de0acbd2e7a3c77a169f7f4f1067934b2a9/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/c3f6fd5a-737b-46c9-98c1-5f206e633e5e.f0db7ecbfc4f7a5ba47fde688a2be3ff12caebd8.zh-cn.xlf", "", "", "aff84ab7-b464-4606-a0ba-ecfa11f6caa9.c65df9e396abea830eb1db4ac21f6ac9e509b63c.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/c520f2bba2caa23f5b9bc262d9c8f73eb53bdab6/e2e/c3f6fd5a-737b-46c9-98c1-5f206e633e5e.md", "", "", "aff84ab7-b464-4606-a0ba-ecfa11f6caa9.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/fda026d82f9baf78f1330c963af240a03e89ea7f/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/c3f6fd5a-737b-46c9-98c1-5f206e633e5e.f0db7ecbfc4f7a5ba47fde688a2be3ff12caebd8.zh-cn.xlf", "", "", "aff84ab7-b464-4606-a0ba-ecfa11f6caa9.c65df9e396abea830eb1db4ac21f6ac9e509b63c.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/11886f77cf3d268b498142352f706872679342af/e2e/ffff4f896d8c-d0bd-4be3-aacb-602751cc319c.md", "", "", "ffffffc9b3ff4e-2f02-46c4-a776-5467936972a9.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3a9e658be81d8a71e56562cc4a5897933818b65e/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/aff84ab7-b464-4606-a0ba-ecfa11f6caa9.c65df9e396abea830eb1db4ac21f6ac9e509b63c.zh-cn.xlf", "", "", "aff84ab7-b464-4606-a0ba-ecfa11f6caa9.c65df9e396abea830eb1db4ac21f6ac9e509b63c.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/fd696c347c514e1bccfb454a68205319abefbe9a/e2e/aff84ab7-b464-4606-a0ba-ecfa11f6caa9.md", "", "", "aff84ab7-b464-4606-a0ba-ecfa11f6caa9.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/2fd4598c2c46cdd4f1ec3aae7e174d3bf8610d8e/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/aff84ab7-b464-4606-a0ba-ecfa11f6caa9.c65df9e396abea830eb1db4ac21f6ac9e509b63c.zh-cn.xlf", "", "", "aff84ab7-b464-4606-a0ba-ecfa11f6caa9.c65df9e396abea830eb1db4ac21f6ac9e509b63c.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/11886f77cf3d268b498142352f706872679342af/e2e/ffffffc9b3ff4e-2f02-46c4-a776-5467936972a9.md", "", "", "c3f6fd5a-737b-46c9-98c1-5f206e633e5e.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3a9e658be81d8a71e56562cc4a5897933818b65e/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/aff84ab7-b464-4606-a0ba-ecfa11f6caa9.c65df9e396abea830eb1db4ac21f6ac9e509b63c.zh-cn.xlf", "", "", "c3f6fd5a-737b-46c9-98c1-5f206e633e5e.f0db7ecbfc4f7a5ba47fde688a2be3ff12caebd8.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("E4"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/fd696c347c514e1bccfb454a68205319abefbe9a/e2e/aff84ab7-b464-4606-a0ba-ecfa11f6caa9.md", "", "", "c3f6fd5a-737b-46c9-98c1-5f206e633e5e.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("F4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/2fd4598c2c46cdd4f1ec3aae7e174d3bf8610d8e/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/aff84ab7-b464-4606-a0ba-ecfa11f6caa9.c65df9e396abea830eb1db4ac21f6ac9e509b63c.zh-cn.xlf", "", "", "c3f6fd5a-737b-46c9-98c1-5f206e633e5e.f0db7ecbfc4f7a5ba47fde688a2be3ff12caebd8.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/11886f77cf3d268b498142352f706872679342af/.localization-config", "", "", ".localization-config") | Out-Null

# ===================== Sheet 3: de-de =====================
$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Range("A1").Hyperlinks.Delete()

$ws3.Range("B2").Value = "Handed back: in sync with en-US"
$ws3.Range("D2").Value = "2016-01-26 12:17:24"
$ws3.Range("G2").Value = "2016-01-26 12:18:21"
$ws3.Range("H2").Value = "Include"
$ws3.Range("B3").Value = "Handed back: in sync with en-US"
$ws3.Range("D3").Value = "2016-01-26 12:17:24"
$ws3.Range("G3").Value = "2016-01-26 12:18:21"
$ws3.Range("H3").Value = "Include"
$ws3.Range("B4").Value = "Ready for handoff"
$ws3.Range("D4").Value = "2016-01-26 12:23:35"
$ws3.Range("G4").Value = "2016-01-26 12:22:29"
$ws3.Range("H4").Value = "Include"
$ws3.Range("B5").Value = "Not to be localized"
$ws3.Range("D5").Value = "0001-01-01 00:00:00"
$ws3.Range("G5").Value = "0001-01-01 00:00:00"
$ws3.Range("H5").Value = "Ignored"

$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/11886f77cf3d268b498142352f706872679342af/e2e/c3f6fd5a-737b-46c9-98c1-5f206e633e5e.md", "", "", "ffff4f896d8c-d0bd-4be3-aacb-602751cc319c.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/831ea6e74b72f08723df48cc4187c967c2f2d8fe/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/c3f6fd5a-737b-46c9-98c1-5f206e633e5e.f0db7ecbfc4f7a5ba47fde688a2be3ff12caebd8.de-de.xlf", "", "", "aff84ab7-b464-4606-a0ba-ecfa11f6caa9.c65df9e396abea830eb1db4ac21f6ac9e509b63c.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/ef56f79981faafd3ad1da53fc558b782208f483c/e2e/c3f6fd5a-737b-46c9-98c1-5f206e633e5e.md", "", "", "aff84ab7-b464-4606-a0ba-ecfa11f6caa9.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/74d6d42b8c18faa488b6b249a8356b9b0bca24c2/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/c3f6fd5a-737b-46c9-98c1-5f206e633e5e.f0db7ecbfc4f7a5ba47fde688a2be3ff12caebd8.de-de.xlf", "", "", "aff84ab7-b464-4606-a0ba-ecfa11f6caa9.c65df9e396abea830eb1db4ac21f6ac9e509b63c.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/11886f77cf3d268b498142352f706872679342af/e2e/ffff4f896d8c-d0bd-4be3-aacb-602751cc319c.md", "", "", "ffffffc9b3ff4e-2f02-46c4-a776-5467936972a9.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/70fbc5520a1830457cc4ae56152d86c9e920b2d6/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/aff84ab7-b464-4606-a0ba-ecfa11f6caa9.c65df9e396abea830eb1db4ac21f6ac9e509b63c.de-de.xlf", "", "", "aff84ab7-b464-4606-a0ba-ecfa11f6caa9.c65df9e396abea830eb1db4ac21f6ac9e509b63c.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/33bdf39446b74e0bec57e4cfc8681470fcfe7dd4/e2e/aff84ab7-b464-4606-a0ba-ecfa11f6caa9.md", "", "", "aff84ab7-b464-4606-a0ba-ecfa11f6caa9.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/c3712a046656506735803cef3fa1cf97f7a02a92/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/aff84ab7-b464-4606-a0ba-ecfa11f6caa9.c65df9e396abea830eb1db4ac21f6ac9e509b63c.de-de.xlf", "", "", "aff84ab7-b464-4606-a0ba-ecfa11f6caa9.c65df9e396abea830eb1db4ac21f6ac9e509b63c.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/11886f77cf3d268b498142352f706872679342af/e2e/ffffffc9b3ff4e-2f02-46c4-a776-5467936972a9.md", "", "", "c3f6fd5a-737b-46c9-98c1-5f206e633e5e.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/70fbc5520a1830457cc4ae56152d86c9e920b2d6/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/aff84ab7-b464-4606-a0ba-ecfa11f6caa9.c65df9e396abea830eb1db4ac21f6ac9e509b63c.de-de.xlf", "", "", "c3f6fd5a-737b-46c9-98c1-5f206e633e5e.f0db7ecbfc4f7a5ba47fde688a2be3ff12caebd8.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("E4"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/33bdf39446b74e0bec57e4cfc8681470fcfe7dd4/e2e/aff84ab7-b464-4606-a0ba-ecfa11f6caa9.md", "", "", "c3f6fd5a-737b-46c9-98c1-5f206e633e5e.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("F4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/c3712a046656506735803cef3fa1cf97f7a02a92/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/aff84ab7-b464-4606-a0ba-ecfa11f6caa9.c65df9e396abea830eb1db4ac21f6ac9e509b63c.de-de.xlf", "", "", "c3f6fd5a-737b-46c9-98c1-5f206e633e5e.f0db7ecbfc4f7a5ba47fde688a2be3ff12caebd8.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/11886f77cf3d268b498142352f706872679342af/.localization-config", "", "", ".localization-config") | Out-Null
